$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'40.020.15"
$ws.Range("E2").Value = "  +2.51%  "

# Row 3
$ws.Range("D3").Value = "'2.230.82"
$ws.Range("E3").Value = "  +1.32%  "

# Row 4
$ws.Range("E4").Value = "  +0.16%  "

# Row 5
$ws.Range("D5").Value = "'293.59"
$ws.Range("E5").Value = "  -0.84%  "

# Row 6
$ws.Range("D6").Value = "'86.69"
$ws.Range("E6").Value = "  +5.10%  "

# Row 7
$ws.Range("D7").Value = "'0.515"
$ws.Range("E7").Value = "  +1.35%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").Value = "'0.472"
$ws.Range("E9").Value = "  +1.56%  "

# Row 10
$ws.Range("D10").Value = "'31.11"
$ws.Range("E10").Value = "  +7.86%  "

# Row 11
$ws.Range("E11").Value = "  +2.44%  "

# Row 12
$ws.Range("D12").Value = "'47.09"
$ws.Range("E12").Value = "  -0.27%  "

# Row 13
$ws.Range("E13").Value = "  +1.28%  "

# Row 14
$ws.Range("D14").Value = "'6.41"
$ws.Range("E14").Value = "  +2.63%  "

# Row 15
$ws.Range("D15").Value = "'2.573.75"
$ws.Range("E15").Value = "  +1.50%  "

# Row 16
$ws.Range("D16").Value = "'14.08"
$ws.Range("E16").Value = "  +0.19%  "

# Row 17
$ws.Range("D17").Value = "'2.244.20"
$ws.Range("E17").Value = "  +2.00%  "

# Row 18
$ws.Range("D18").Value = "'0.730"
$ws.Range("E18").Value = "  +2.66%  "

# Row 19
$ws.Range("D19").Value = "'39.934.38"
$ws.Range("E19").Value = "  +2.51%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0891"
$ws.Range("E20").Value = "  +2.73%  "

# Row 21
$ws.Range("D21").Value = "'11.12"
$ws.Range("E21").Value = "  +8.90%  "

# Row 22
$ws.Range("D22").Value = "'5.82"
$ws.Range("E22").Value = "  +2.27%  "

# Row 23
$ws.Range("D23").Value = "'65.44"
$ws.Range("E23").Value = "  +1.36%  "

# Row 24
$ws.Range("D24").Value = "'235.59"
$ws.Range("E24").Value = "  +3.70%  "

# Row 25
$ws.Range("E25").Value = "  -0.11%  "

# Row 26
$ws.Range("E26").Value = "  +2.89%  "

# Row 27
$ws.Range("E27").Value = "  +4.06%  "

# Row 28
$ws.Range("D28").Value = "'22.81"
$ws.Range("E28").Value = "  +1.61%  "

# Row 29
$ws.Range("E29").Value = "  +5.60%  "

# Row 30
$ws.Range("E30").Value = "  +3.19%  "

# Row 31
$ws.Range("D31").Value = "'33.21"
$ws.Range("E31").Value = "  +4.26%  "

# Row 32
$ws.Range("D32").Value = "'152.40"
$ws.Range("E32").Value = "  +1.99%  "

# Row 33
$ws.Range("E33").Value = "  +0.02%  "

# Row 34
$ws.Range("D34").Value = "'4.93"
$ws.Range("E34").Value = "  +2.28%  "

# Row 35
$ws.Range("D35").Value = "'0.0720"
$ws.Range("E35").Value = "  +3.93%  "

# Row 36
$ws.Range("D36").Value = "'2.38"
$ws.Range("E36").Value = "  +2.93%  "

# Row 37
$ws.Range("E37").Value = "  +9.66%  "

# Row 38
$ws.Range("D38").Value = "'2.82"
$ws.Range("E38").Value = "  +7.28%  "

# Row 39
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.112"
$ws.Range("E39").Value = "  +2.58%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.100"
$ws.Range("E40").Value = "  +4.86%  "

# Row 41
$ws.Range("E41").Value = "  +5.96%  "

# Row 42
$ws.Range("D42").Value = "'3.81"
$ws.Range("E42").Value = "  +3.66%  "

# Row 43
$ws.Range("D43").Value = "'2.050.67"
$ws.Range("E43").Value = "  +7.50%  "

# Row 44
$ws.Range("D44").Value = "'2.13"
$ws.Range("E44").Value = "  +5.50%  "

# Row 45
$ws.Range("D45").Value = "'0.0269"
$ws.Range("E45").Value = "  +4.42%  "

# Row 46
$ws.Range("D46").Value = "'9.85"
$ws.Range("E46").Value = "  +11.06%  "

# Row 47
$ws.Range("D47").Value = "'17.18"
$ws.Range("E47").Value = "  +8.74%  "

# Row 48
$ws.Range("D48").Value = "'2.60"
$ws.Range("E48").Value = "  -0.79%  "

# Row 49
$ws.Range("D49").Value = "'2.436.43"
$ws.Range("E49").Value = "  +1.01%  "

# Row 50
$ws.Range("D50").Value = "'71.82"
$ws.Range("E50").Value = "  +2.56%  "

# Row 51
$ws.Range("D51").Value = "'89.25"
$ws.Range("E51").Value = "  +3.34%  "
